$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions bot).
# Columns B (Coin), C (Link), D (Price) and E (Volume(1h)) are plain
# text cells in this sheet. Several Price values look like numbers
# (e.g. "1.00", "115.68") so we briefly switch those cells to Text
# number format while assigning the value (preventing Excel from
# re-interpreting them as numeric values / dropping trailing zeros),
# then restore the default "Normal" style so the cell formatting
# is left exactly as it was before the edit.

$ws.Range("D2").Value = "52.182.92"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "2.801.25"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.548"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0863"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("D15").Value = "3.249.63"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").Value = "2.803.33"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.891"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "52.161.13"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.83%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  +28.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0833"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  +10.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.54%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "2.055.54"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.958"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
